# Add a new "Save" column (H) to the sheet, mirroring the style of the
# other header cells and filling in the per-row save flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the rest of row 1 (bold/bordered/centered header style).
# Copy the formatting from the neighboring header cell (G1) so the new
# header cell reuses the same cell style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Per-row "Save" flag values (0/1) for rows 2-15.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
